# Apply the "Add files via upload" edit to slide 1 of the presentation:
#  1. Fix the title typo "RECARP" -> "RECAP".
#  2. Change the subtitle text to "By: Group 2" and push the trailing
#     (formerly shared) paragraph mark onto its own empty paragraph.
#  3. Add a new rounded-rectangle accent shape to the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Title fix -----------------------------------------------------
$title = $s.Shapes.Item("Title 1")
$title.TextFrame.TextRange.Text = "RECAP for the virtual seminar series"

# --- 2. Subtitle text + paragraph split --------------------------------
$subtitle = $s.Shapes.Item("Subtitle 2")
$subtitleRange = $subtitle.TextFrame.TextRange
$subtitleRange.Text = "By: Group 2"
$null = $subtitleRange.InsertAfter([char]13)

# --- 3. New rounded-rectangle shape -------------------------------------
$EMUS_PER_POINT = 12700.0
$msoShapeRoundedRectangle = 5

$offX = 242047 / $EMUS_PER_POINT
$offY = 0 / $EMUS_PER_POINT
$extCx = 2622177 / $EMUS_PER_POINT
$extCy = 685800 / $EMUS_PER_POINT

$rect = $s.Shapes.AddShape($msoShapeRoundedRectangle, $offX, $offY, $extCx, $extCy)
$rect.Name = "Rectangle: Rounded Corners 1"
$rect.TextFrame.TextRange.ParagraphFormat.Alignment = 2   # ppAlignCenter
$rect.TextFrame.VerticalAnchor = 3                          # msoAnchorMiddle
